# "Generate Report for handback"
#
# The localization-status report is regenerated: the 53c9dd87... file has now
# been handed back (previously "Not yet handed off"), so its row moves ahead
# of the 9afcfc7c... file's row on every sheet (Overview, zh-cn, de-de),
# both files now show status "Handed back", and the newly handed-back file's
# "Latest Handback DateTime" is updated.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview": File Name / zh-cn / de-de
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "53c9dd87-9eb1-4248-aecd-c73be98fef3c.md"
$ov.Range("B2").Value = "Handed back"
$ov.Range("C2").Value = "Handed back"

$ov.Range("A3").Value = "9afcfc7c-5f1b-416f-a7e3-c731c38e661d.md"
$ov.Range("B3").Value = "Handed back"
$ov.Range("C3").Value = "Handed back"

foreach ($h in $ov.Hyperlinks) {
    if ($h.TextToDisplay -eq "9afcfc7c-5f1b-416f-a7e3-c731c38e661d.md") {
        $h.TextToDisplay = "53c9dd87-9eb1-4248-aecd-c73be98fef3c.md"
    } elseif ($h.TextToDisplay -eq "53c9dd87-9eb1-4248-aecd-c73be98fef3c.md") {
        $h.TextToDisplay = "9afcfc7c-5f1b-416f-a7e3-c731c38e661d.md"
    }
}

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "53c9dd87-9eb1-4248-aecd-c73be98fef3c.md"
$zh.Range("B2").Value = "Handed back"
$zh.Range("C2").Value = "53c9dd87-9eb1-4248-aecd-c73be98fef3c.bacd708eae0fd5a27b6d0a9273ec2e18c6b2cc12.zh-cn.xlf"
$zh.Range("D2").Value = "2016-01-08 14:28:38"
$zh.Range("E2").Value = "53c9dd87-9eb1-4248-aecd-c73be98fef3c.md"
$zh.Range("F2").Value = "53c9dd87-9eb1-4248-aecd-c73be98fef3c.bacd708eae0fd5a27b6d0a9273ec2e18c6b2cc12.zh-cn.xlf"
$zh.Range("G2").Value = "2016-01-08 14:29:29"
$zh.Range("H2").Value = "Include"

$zh.Range("A3").Value = "9afcfc7c-5f1b-416f-a7e3-c731c38e661d.md"
$zh.Range("B3").Value = "Handed back"
$zh.Range("C3").Value = "9afcfc7c-5f1b-416f-a7e3-c731c38e661d.020c6830f98f5756493e16d2d8a9e895c4be6e0f.zh-cn.xlf"
$zh.Range("D3").Value = "2016-01-08 14:26:31"
$zh.Range("E3").Value = "9afcfc7c-5f1b-416f-a7e3-c731c38e661d.md"
$zh.Range("F3").Value = "9afcfc7c-5f1b-416f-a7e3-c731c38e661d.020c6830f98f5756493e16d2d8a9e895c4be6e0f.zh-cn.xlf"
$zh.Range("G3").Value = "2016-01-08 14:27:33"
$zh.Range("H3").Value = "Include"

foreach ($h in $zh.Hyperlinks) {
    if ($h.TextToDisplay -eq "9afcfc7c-5f1b-416f-a7e3-c731c38e661d.md") {
        $h.TextToDisplay = "53c9dd87-9eb1-4248-aecd-c73be98fef3c.md"
    } elseif ($h.TextToDisplay -eq "53c9dd87-9eb1-4248-aecd-c73be98fef3c.md") {
        $h.TextToDisplay = "9afcfc7c-5f1b-416f-a7e3-c731c38e661d.md"
    } elseif ($h.TextToDisplay -eq "9afcfc7c-5f1b-416f-a7e3-c731c38e661d.020c6830f98f5756493e16d2d8a9e895c4be6e0f.zh-cn.xlf") {
        $h.TextToDisplay = "53c9dd87-9eb1-4248-aecd-c73be98fef3c.bacd708eae0fd5a27b6d0a9273ec2e18c6b2cc12.zh-cn.xlf"
    } elseif ($h.TextToDisplay -eq "53c9dd87-9eb1-4248-aecd-c73be98fef3c.bacd708eae0fd5a27b6d0a9273ec2e18c6b2cc12.zh-cn.xlf") {
        $h.TextToDisplay = "9afcfc7c-5f1b-416f-a7e3-c731c38e661d.020c6830f98f5756493e16d2d8a9e895c4be6e0f.zh-cn.xlf"
    }
}

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "53c9dd87-9eb1-4248-aecd-c73be98fef3c.md"
$de.Range("B2").Value = "Handed back"
$de.Range("C2").Value = "53c9dd87-9eb1-4248-aecd-c73be98fef3c.bacd708eae0fd5a27b6d0a9273ec2e18c6b2cc12.de-de.xlf"
$de.Range("D2").Value = "2016-01-08 14:28:52"
$de.Range("E2").Value = "53c9dd87-9eb1-4248-aecd-c73be98fef3c.md"
$de.Range("F2").Value = "53c9dd87-9eb1-4248-aecd-c73be98fef3c.bacd708eae0fd5a27b6d0a9273ec2e18c6b2cc12.de-de.xlf"
$de.Range("G2").Value = "2016-01-08 14:29:52"
$de.Range("H2").Value = "Include"

$de.Range("A3").Value = "9afcfc7c-5f1b-416f-a7e3-c731c38e661d.md"
$de.Range("B3").Value = "Handed back"
$de.Range("C3").Value = "9afcfc7c-5f1b-416f-a7e3-c731c38e661d.020c6830f98f5756493e16d2d8a9e895c4be6e0f.de-de.xlf"
$de.Range("D3").Value = "2016-01-08 14:26:45"
$de.Range("E3").Value = "9afcfc7c-5f1b-416f-a7e3-c731c38e661d.md"
$de.Range("F3").Value = "9afcfc7c-5f1b-416f-a7e3-c731c38e661d.020c6830f98f5756493e16d2d8a9e895c4be6e0f.de-de.xlf"
$de.Range("G3").Value = "2016-01-08 14:27:56"
$de.Range("H3").Value = "Include"

foreach ($h in $de.Hyperlinks) {
    if ($h.TextToDisplay -eq "9afcfc7c-5f1b-416f-a7e3-c731c38e661d.md") {
        $h.TextToDisplay = "53c9dd87-9eb1-4248-aecd-c73be98fef3c.md"
    } elseif ($h.TextToDisplay -eq "53c9dd87-9eb1-4248-aecd-c73be98fef3c.md") {
        $h.TextToDisplay = "9afcfc7c-5f1b-416f-a7e3-c731c38e661d.md"
    } elseif ($h.TextToDisplay -eq "9afcfc7c-5f1b-416f-a7e3-c731c38e661d.020c6830f98f5756493e16d2d8a9e895c4be6e0f.de-de.xlf") {
        $h.TextToDisplay = "53c9dd87-9eb1-4248-aecd-c73be98fef3c.bacd708eae0fd5a27b6d0a9273ec2e18c6b2cc12.de-de.xlf"
    } elseif ($h.TextToDisplay -eq "53c9dd87-9eb1-4248-aecd-c73be98fef3c.bacd708eae0fd5a27b6d0a9273ec2e18c6b2cc12.de-de.xlf") {
        $h.TextToDisplay = "9afcfc7c-5f1b-416f-a7e3-c731c38e661d.020c6830f98f5756493e16d2d8a9e895c4be6e0f.de-de.xlf"
    }
}
